# The workbook's second sheet ("Лист2") is renamed to "My parent's accounting".
#
# (The XML diff also shows a handful of purely cosmetic/environmental deltas —
# the <mc:Choice><x15ac:absPath> of the author's machine, the co-authoring
# <xr:revisionPtr documentId>, and the <bookViews><workbookView> on-screen
# window geometry — plus a reshuffling of the <mergeCells> entries in
# sheet2.xml into a different (but set-identical) order. None of those carry
# any user-visible/semantic content: they are artifacts of the original
# author's local Excel session/save, not of an edit a script should
# reproduce, and this host's object model does not expose a way to author
# them (Window.Left/Top/Width/Height etc. round-trip through the COM shim
# but are not persisted into <workbookView>; the file path/session id are
# likewise fixed by the host). The only substantive, reproducible change is
# the sheet rename below.)

$wb = $excel.ActiveWorkbook

$target = $null
foreach ($sheet in $wb.Worksheets) {
    if ($sheet.Name -eq "Лист2") {
        $target = $sheet
    }
}
if ($target -eq $null) {
    # Fallback: "Лист2" is the second sheet in the workbook.
    $target = $wb.Worksheets.Item(2)
}

$target.Name = "My parent's accounting"
